# Update cryptos list figures (price & 1h volume change) and fix the
# ordering/data for Aptos vs FirstDigitalUSD (rows 33 and 34).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.824.79"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.358.08"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.39"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.20"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  -2.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.370.60"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0970"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.80"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.775.41"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "55.760.16"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.49"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.384.16"
$ws.Range("E18").Value = "  +24.72%  "
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "310.32"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.23"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.32"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.13"
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.55"
$ws.Range("E29").Value = "  -2.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0706"
$ws.Range("E30").Value = "  -2.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.64"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.73"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.995"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.66"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.838"
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("E39").Value = "  -4.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.18"
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.86"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.92"
$ws.Range("E44").Value = "  -4.97%  "
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0896"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "240.13"
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0478"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.84"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.63"
$ws.Range("E51").Value = "  -2.93%  "
